# 3DES projetos - Adicionadas aulas de 06 a 10
# Extends the FREQ attendance sheet with the third week of classes
# (columns L:P), mirroring the existing PROJ/PROJ/PROJ/PDMO/RMST weekly
# header pattern and the row-2 date strip, then records attendance
# (P = presente, F = falta) for the two days that already have data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FREQ")

# --- Row 1: subject headers for the new week (PROJ, PROJ, PROJ, PDMO, RMST) ---
$ws1.Range("L1").Value = "PROJ"
$ws1.Range("M1").Value = "PROJ"
$ws1.Range("N1").Value = "PROJ"
$ws1.Range("O1").Value = "PDMO"
$ws1.Range("P1").Value = "RMST"

# --- Row 2: class dates for the new week, copying K2's date format/style ---
$ws1.Range("K2").Copy() | Out-Null
$ws1.Range("L2:P2").PasteSpecial(-4122) | Out-Null
$ws1.Range("L2").Value = 44214
$ws1.Range("M2").Value = 44215
$ws1.Range("N2").Value = 44216
$ws1.Range("O2").Value = 44217
$ws1.Range("P2").Value = 44218
$excel.CutCopyMode = 0

# --- Rows 3-20: attendance (P/F) recorded so far for the two new classes ---
$attendance = @{
  3  = @("P", "P")
  4  = @("P", "P")
  5  = @("F", "P")
  6  = @("P", "P")
  7  = @("F", "F")
  8  = @("F", "P")
  9  = @("P", "P")
  10 = @("P", "P")
  11 = @("P", "P")
  12 = @("P", "P")
  13 = @("F", "P")
  14 = @("P", "P")
  15 = @("P", "P")
  16 = @("P", "P")
  17 = @("P", "P")
  18 = @("P", "P")
  19 = @("P", "P")
  20 = @("P", "P")
}

foreach ($r in 3..20) {
  $vals = $attendance[$r]
  $ws1.Range("L$r").Value = $vals[0]
  $ws1.Range("M$r").Value = $vals[1]
}

# --- Selection / active-tab bookkeeping: FREQ becomes the active sheet ---
$ws1.Activate() | Out-Null
$ws1.Range("L20").Select() | Out-Null
